$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values are stored as text in the source sheet, so force
# text entry (matches the existing inline-string / text cell type) instead of
# letting Excel auto-convert numeric-looking strings into real numbers, which
# would silently drop meaningful trailing zeros (e.g. "54.90" -> 54.9).
$priceUpdates = @{
    2 = "69.917.14"
    3 = "3.565.84"
    4 = "0.999"
    5 = "575.46"
    6 = "185.87"
    7 = "3.561.40"
    8 = "0.619"
    10 = "0.184"
    11 = "0.648"
    12 = "54.90"
    13 = "0.0000302"
    14 = "9.50"
    15 = "4.137.45"
    16 = "19.57"
    17 = "3.552.64"
    18 = "69.847.49"
    19 = "12.55"
    21 = "1.03"
    22 = "496.14"
    23 = "19.16"
    24 = "4.90"
    25 = "4.38"
    26 = "95.26"
    27 = "11.33"
    28 = "2.94"
    29 = "9.28"
    30 = "31.58"
    31 = "7.54"
    32 = "66.81"
    33 = "12.03"
    34 = "0.115"
    35 = "566.11"
    36 = "3.14"
    37 = "38.57"
    39 = "0.0₃0789"
    40 = "0.393"
    41 = "3.52"
    42 = "3.18"
    43 = "0.134"
    44 = "2.99"
    45 = "3.225.51"
    46 = "3.46"
    47 = "0.0439"
    48 = "9.55"
    50 = "0.998"
    51 = "3.14"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Volume(1h) column (E): percentage strings already contain non-numeric
# characters (spaces, "%"), so a plain assignment keeps them as text.
$volumeUpdates = @{
    2 = "  -1.54%  "
    3 = "  -2.52%  "
    4 = "  -0.10%  "
    5 = "  -3.52%  "
    6 = "  -4.54%  "
    7 = "  -2.46%  "
    8 = "  -4.43%  "
    9 = "  +0.03%  "
    10 = "  -0.46%  "
    11 = "  -4.07%  "
    12 = "  -6.05%  "
    13 = "  +2.45%  "
    14 = "  -4.80%  "
    15 = "  -2.46%  "
    16 = "  -2.83%  "
    17 = "  -2.93%  "
    18 = "  -1.65%  "
    19 = "  -2.17%  "
    20 = "  -0.95%  "
    21 = "  -3.71%  "
    22 = "  +1.36%  "
    23 = "  +0.74%  "
    24 = "  -7.37%  "
    25 = "  -2.94%  "
    26 = "  +4.12%  "
    27 = "  -1.16%  "
    28 = "  -6.98%  "
    29 = "  -3.61%  "
    30 = "  -4.04%  "
    31 = "  -3.52%  "
    32 = "  +0.56%  "
    33 = "  -2.22%  "
    34 = "  -6.28%  "
    35 = "  -10.12%  "
    36 = "  +11.19%  "
    37 = "  -4.44%  "
    38 = "  -0.03%  "
    39 = "  -5.23%  "
    40 = "  -4.92%  "
    41 = "  -2.19%  "
    42 = "  +3.57%  "
    43 = "  -10.14%  "
    44 = "  -5.50%  "
    46 = "  +4.16%  "
    47 = "  -3.73%  "
    48 = "  +0.82%  "
    49 = "  -3.20%  "
    50 = "  -0.26%  "
    51 = "  -3.85%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

